$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Formula = "'245.89"
$ws.Range("G2").Formula = "'11"

# Row 3
$ws.Range("D3").Formula = "'23.93"
$ws.Range("G3").Formula = "'11"

# Row 4
$ws.Range("D4").Formula = "'5.278"
$ws.Range("G4").Formula = "'11"

# Row 5
$ws.Range("D5").Formula = "'0.05736"
$ws.Range("G5").Formula = "'11"

# Row 6
$ws.Range("D6").Formula = "'6.492"
$ws.Range("G6").Formula = "'11"

# Row 7
$ws.Range("D7").Formula = "'3.153"
$ws.Range("G7").Formula = "'11"

# Row 8
$ws.Range("D8").Formula = "'0.8158"
$ws.Range("G8").Formula = "'11"

# Row 9
$ws.Range("D9").Formula = "'0.8563"
$ws.Range("G9").Formula = "'11"

# Row 10
$ws.Range("G10").Formula = "'11"

# Row 11
$ws.Range("D11").Formula = "'0.07000"
$ws.Range("G11").Formula = "'11"

# Row 12
$ws.Range("D12").Formula = "'0.03203"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("G12").Formula = "'11"

# Row 13
$ws.Range("D13").Formula = "'0.02871"
$ws.Range("G13").Formula = "'11"

# Row 14
$ws.Range("D14").Formula = "'0.09398"
$ws.Range("G14").Formula = "'11"

# Row 15
$ws.Range("D15").Formula = "'3.823"
$ws.Range("G15").Formula = "'11"

# Row 16
$ws.Range("D16").Formula = "'0.001529"
$ws.Range("G16").Formula = "'11"

# Row 17
$ws.Range("D17").Formula = "'0.04690"
$ws.Range("G17").Formula = "'11"

# Row 18
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D18").Formula = "'0.006233"
$ws.Range("E18").Value = "17TigerCashTCH"
$ws.Range("G18").Formula = "'11"

# Row 19
$ws.Range("B19").Value = "BitKan"
$ws.Range("C19").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D19").Formula = "'0.001244"
$ws.Range("E19").Value = "18BitKanKAN"
$ws.Range("G19").Formula = "'11"

# Row 20
$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D20").Formula = "'0.004783"
$ws.Range("E20").Value = "19HotbitTokenHTB"
$ws.Range("G20").Formula = "'11"

# Row 21
$ws.Range("B21").Value = "NitroEx"
$ws.Range("C21").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D21").Formula = "'0.00006998"
$ws.Range("E21").Value = "20NitroExNTX"
$ws.Range("G21").Formula = "'11"

# Row 22
$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D22").Formula = "'3.532"
$ws.Range("E22").Value = "21LEOLEO"
$ws.Range("G22").Formula = "'11"

# Row 23
$ws.Range("B23").Value = "BTSEToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D23").Formula = "'2.149"
$ws.Range("E23").Value = "22BTSETokenBTSE"
$ws.Range("G23").Formula = "'11"

# Row 24
$ws.Range("B24").Value = "One"
$ws.Range("C24").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D24").Formula = "'0.0005988"
$ws.Range("E24").Value = "23OneONE"
$ws.Range("G24").Formula = "'11"

# Row 25
$ws.Range("D25").Formula = "'0.3175"
$ws.Range("G25").Formula = "'11"

# Row 26
$ws.Range("G26").Formula = "'11"

# Row 27
$ws.Range("G27").Formula = "'11"

# Row 28
$ws.Range("G28").Formula = "'11"

# Row 29
$ws.Range("G29").Formula = "'11"

# Row 30
$ws.Range("G30").Formula = "'11"

# Row 31
$ws.Range("G31").Formula = "'11"

# Row 32
$ws.Range("G32").Formula = "'11"

# Row 33
$ws.Range("G33").Formula = "'11"

# Row 34
$ws.Range("G34").Formula = "'11"

# Row 35
$ws.Range("G35").Formula = "'11"

# Row 36
$ws.Range("G36").Formula = "'11"

# Row 37
$ws.Range("G37").Formula = "'11"

# Row 38
$ws.Range("G38").Formula = "'11"

# Row 39
$ws.Range("G39").Formula = "'11"

# Row 40
$ws.Range("D40").Formula = "'0.03705"
$ws.Range("G40").Formula = "'11"

# Row 41
$ws.Range("D41").Formula = "'0.006338"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
$ws.Range("G41").Formula = "'11"

# Row 42
$ws.Range("D42").Formula = "'0.1054"
$ws.Range("G42").Formula = "'11"

# Row 43
$ws.Range("D43").Formula = "'0.002211"
$ws.Range("G43").Formula = "'11"

# Row 44
$ws.Range("D44").Formula = "'0.008709"
$ws.Range("G44").Formula = "'11"

# Row 45
$ws.Range("D45").Formula = "'0.00005488"
$ws.Range("G45").Formula = "'11"

# Row 46
$ws.Range("G46").Formula = "'11"

# Row 47
$ws.Range("D47").Formula = "'0.3883"
$ws.Range("G47").Formula = "'11"

# Row 48
$ws.Range("D48").Formula = "'0.002572"
$ws.Range("G48").Formula = "'11"

# Row 49
$ws.Range("G49").Formula = "'11"

# Row 50
$ws.Range("G50").Formula = "'11"

# Row 51
$ws.Range("G51").Formula = "'11"
